# "Add row in excel sheet"
# Adds new login-credential rows and a drop-down source list of Data
# Structures topics to Sheet1, drops the (now stale) mailto hyperlinks on
# that sheet, and leaves Sheet1 as the active sheet/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New credential rows
$ws.Range("A11").Value = "milestonemavericks"
$ws.Range("B11").Value = "Welcome@1"

$ws.Range("A12").Value = "mmuser1"
$ws.Range("B12").Value = "Welcome@123"

# Row 13 intentionally left blank as a separator.

# Drop-down source list (Data Structures topics)
$ws.Range("A14").Value = "Drop-down"
$ws.Range("A15").Value = "Data Structures"
$ws.Range("A16").Value = "Arrays"
$ws.Range("A17").Value = "Linked List"
$ws.Range("A18").Value = "Stack"
$ws.Range("A19").Value = "Queue"
$ws.Range("A20").Value = "Tree"
$ws.Range("A21").Value = "Graph"

# Hyperlinks on this sheet are no longer needed - drop them.
$ws.Hyperlinks.Delete() | Out-Null

# Make Sheet1 the active sheet/selection.
$ws.Activate() | Out-Null
$ws.Range("C17").Select() | Out-Null
